$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.409.92'
Set-TextValue 'E2' '  +0.06%  '
Set-TextValue 'D3' '1.849.94'
Set-TextValue 'E3' '  +0.15%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '240.50'
Set-TextValue 'E5' '  +0.01%  '
Set-TextValue 'D6' '0.6294'
Set-TextValue 'E6' '  +0.15%  '
Set-TextValue 'D8' '0.07643'
Set-TextValue 'E8' '  +0.51%  '
Set-TextValue 'D9' '0.2911'
Set-TextValue 'E9' '  -0.60%  '
Set-TextValue 'D10' '24.83'
Set-TextValue 'E10' '  +1.48%  '
Set-TextValue 'D11' '2.160.62'
Set-TextValue 'E11' '  +16.79%  '
Set-TextValue 'E12' '  +0.08%  '
Set-TextValue 'D13' '5.033'
Set-TextValue 'E13' '  +0.63%  '
Set-TextValue 'D14' '0.6811'
Set-TextValue 'E14' '  +0.35%  '
Set-TextValue 'D15' '0.00001075'
Set-TextValue 'E15' '  -0.75%  '
Set-TextValue 'D16' '83.33'
Set-TextValue 'E16' '  -0.49%  '
Set-TextValue 'D17' '6.169'
Set-TextValue 'E17' '  -0.09%  '
Set-TextValue 'D18' '29.512.28'
Set-TextValue 'E18' '  +0.35%  '
Set-TextValue 'D19' '228.30'
Set-TextValue 'E19' '  -0.07%  '
Set-TextValue 'D20' '12.33'
Set-TextValue 'E20' '  -0.90%  '
Set-TextValue 'E21' '  +0.06%  '
Set-TextValue 'E22' '  -0.03%  '
Set-TextValue 'E23' '  +0.02%  '
Set-TextValue 'D24' '158.00'
Set-TextValue 'E24' '  +0.46%  '
Set-TextValue 'D25' '0.1383'
Set-TextValue 'D26' '8.429'
Set-TextValue 'E26' '  +0.93%  '
Set-TextValue 'D27' '17.67'
Set-TextValue 'E27' '  +0.26%  '
Set-TextValue 'D28' '1.376'
Set-TextValue 'E28' '  +6.14%  '
Set-TextValue 'D29' '1.461'
Set-TextValue 'E29' '  -0.25%  '
Set-TextValue 'D30' '0.05612'
Set-TextValue 'E30' '  +0.53%  '
Set-TextValue 'D31' '4.132'
Set-TextValue 'E31' '  +0.82%  '
Set-TextValue 'D32' '4.057'
Set-TextValue 'E32' '  +0.68%  '
Set-TextValue 'D33' '1.842'
Set-TextValue 'E33' '  -0.23%  '
Set-TextValue 'D34' '1.164'
Set-TextValue 'E34' '  +0.75%  '
Set-TextValue 'D35' '0.6930'
Set-TextValue 'E35' '  -2.29%  '
Set-TextValue 'D36' '2.589'
Set-TextValue 'E36' '  +0.14%  '
Set-TextValue 'D37' '0.01804'
Set-TextValue 'E37' '  +0.25%  '
Set-TextValue 'D38' '1.227.63'
Set-TextValue 'E38' '  -0.44%  '
Set-TextValue 'D39' '2.718'
Set-TextValue 'E39' '  -1.92%  '
Set-TextValue 'D40' '6.449'
Set-TextValue 'E40' '  +0.61%  '
Set-TextValue 'D41' '0.9068'
Set-TextValue 'E41' '  -0.02%  '
Set-TextValue 'E42' '  +0.08%  '
Set-TextValue 'D43' '101.73'
Set-TextValue 'E43' '  -0.05%  '
Set-TextValue 'D44' '65.98'
Set-TextValue 'E44' '  -0.07%  '
Set-TextValue 'E45' '  -0.72%  '
Set-TextValue 'D46' '7.192'
Set-TextValue 'D47' '0.4021'
Set-TextValue 'E47' '  +0.09%  '
Set-TextValue 'D48' '0.1153'
Set-TextValue 'E48' '  +2.91%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '8.988'
Set-TextValue 'E49' '  -0.34%  '
Set-TextValue 'B50' 'RenderToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '1.682'
Set-TextValue 'E50' '  +0.35%  '
Set-TextValue 'D51' '0.05703'
Set-TextValue 'E51' '  -0.11%  '
